# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# before column N, pushing the existing "Late" / "heading" / "Outstanding"
# columns one slot to the right (N->O, O->P, P->Q). The new column inherits
# its width from its left neighbour (column M, width 11).
#
# The "Repayment schedule" sheet also becomes the active/selected sheet
# (it was previously "Edit Repayment Schedule"), with the active cell
# moved to J13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Repayment schedule")

# Insert a new blank column at N; existing N/O/P data+styles shift to O/P/Q.
$ws.Columns("N").Insert() | Out-Null

# Match the inherited column width (stored width "11" in the xlsx, which is
# ColumnWidth 11 minus Excel's fixed ~5/6 character padding offset).
$ws.Columns("N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab, with J13 selected
# (previously "Edit Repayment Schedule" held the selected tab).
$ws.Activate() | Out-Null
$ws.Range("J13").Select() | Out-Null
